# feat: add 2022-Q3 data
#
# 1. Update the "总计" (total) summary sheet: a new row for 2022-Q3 is
#    inserted at the top of the data block (row 2), pushing the other
#    quarters' summary figures down by one row. A brand-new trailing row
#    (2020-Q4, index 5) is appended at the bottom.
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计",
#    containing the quarter's per-fund holdings detail, and shift the
#    remaining quarter sheets along (handled automatically by Excel).

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing a genuine text/string
# cell type (matches the source data, where figures such as "25.09" or
# fund codes like "001556" are stored as text, not numbers), without
# leaving a stray NumberFormat behind on the cell.
function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = [string]$val
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. "总计" sheet — shift existing quarters down one row and splice in
#    the new 2022-Q3 summary numbers at the top.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryRows = @(
    @{ A = 0; B = "2022-Q3"; C = 6;  D = 0.89 },
    @{ A = 1; B = "2022-Q2"; C = 2;  D = 0.39 },
    @{ A = 2; B = "2022-Q1"; C = 13; D = 0.34 },
    @{ A = 3; B = "2021-Q4"; C = 2;  D = 0.11 },
    @{ A = 4; B = "2021-Q3"; C = 4;  D = 0.12 },
    @{ A = 5; B = "2020-Q4"; C = 1;  D = 0.24 }
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $summary $r 2 $row.B
    $summary.Cells.Item($r, 3).Value = $row.C
    $summary.Cells.Item($r, 4).Value = $row.D
}

# Row 7 (index 5, "2020-Q4") is brand new - give its A cell the same
# look (bold/border/center) as the other index cells above it (A2:A6).
$summary.Cells.Item(2, 1).Copy()
$summary.Cells.Item(7, 1).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Brand-new "2022-Q3" sheet, positioned right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    Set-TextCell $q3 1 $col $headers[$col - 2]
}

# Match the look of the header row / index column used throughout the
# other quarter sheets: bold text, thin box border, centered top-aligned.
$headerRange = $q3.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$detailRows = @(
    @{ A = 0; B = "001556"; C = "天弘中证500指数增强A";             D = "25.09"; E = "94.15"; F = "1.56"; G = "0.3914"; H = 6 },
    @{ A = 1; B = "001557"; C = "天弘中证500指数增强C";             D = "12.94"; E = "94.15"; F = "1.56"; G = "0.2019"; H = 6 },
    @{ A = 2; B = "005994"; C = "国投瑞银中证500指数量化增强A"; D = "13.36"; E = "88.67"; F = "1.13"; G = "0.1510"; H = 10 },
    @{ A = 3; B = "005396"; C = "中金丰硕混合";                     D = "1.77";  E = "76.61"; F = "3.99"; G = "0.0706"; H = 10 },
    @{ A = 4; B = "007089"; C = "国投瑞银中证500指数量化增强C"; D = "4.45";  E = "88.67"; F = "1.13"; G = "0.0503"; H = 10 },
    @{ A = 5; B = "000270"; C = "建信灵活配置混合";                 D = "2.27";  E = "94.21"; F = "0.90"; G = "0.0204"; H = 8 }
)

for ($i = 0; $i -lt $detailRows.Count; $i++) {
    $r = $i + 2
    $row = $detailRows[$i]
    $q3.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $q3 $r 2 $row.B
    Set-TextCell $q3 $r 3 $row.C
    Set-TextCell $q3 $r 4 $row.D
    Set-TextCell $q3 $r 5 $row.E
    Set-TextCell $q3 $r 6 $row.F
    Set-TextCell $q3 $r 7 $row.G
    $q3.Cells.Item($r, 8).Value = $row.H
}

$indexRange = $q3.Range("A2:A7")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108  # xlCenter
$indexRange.VerticalAlignment = -4160    # xlTop
$indexRange.Borders.LineStyle = 1
$indexRange.Borders.Weight = 2
